# Weekly update: insert two new Ciruela price rows at the top of the
# data block (rows 69-70), pushing the existing rows 69-148 down to 150.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 69 (this shifts the old rows 69..148 down
# to 71..150, carrying all their original values/formatting with them).
$ws.Rows.Item(69).Insert()
$ws.Rows.Item(69).Insert()

# --- Fill in the new row 69 ---
$ws.Cells.Item(69, 1).Value  = 9
$ws.Cells.Item(69, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(69, 3).Value  = "Metropolitana"
$ws.Cells.Item(69, 4).Value  = 44942
$ws.Cells.Item(69, 5).Value  = 13
$ws.Cells.Item(69, 6).Value  = "Fruta"
$ws.Cells.Item(69, 7).Value  = 100103
$ws.Cells.Item(69, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(69, 9).Value  = 100103002
$ws.Cells.Item(69, 10).Value = "Ciruela"
$ws.Cells.Item(69, 11).Value = "Black Amber"
$ws.Cells.Item(69, 12).Value = "Especial"
$ws.Cells.Item(69, 13).Value = 280
$ws.Cells.Item(69, 14).Value = 9000
$ws.Cells.Item(69, 15).Value = 9000
$ws.Cells.Item(69, 16).Value = 9000
$ws.Cells.Item(69, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(69, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(69, 19).Value = 900
$ws.Cells.Item(69, 20).Value = 10

# --- Fill in the new row 70 ---
$ws.Cells.Item(70, 1).Value  = 9
$ws.Cells.Item(70, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(70, 3).Value  = "Metropolitana"
$ws.Cells.Item(70, 4).Value  = 44942
$ws.Cells.Item(70, 5).Value  = 13
$ws.Cells.Item(70, 6).Value  = "Fruta"
$ws.Cells.Item(70, 7).Value  = 100103
$ws.Cells.Item(70, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(70, 9).Value  = 100103002
$ws.Cells.Item(70, 10).Value = "Ciruela"
$ws.Cells.Item(70, 11).Value = "Black Amber"
$ws.Cells.Item(70, 12).Value = "Primera"
$ws.Cells.Item(70, 13).Value = 220
$ws.Cells.Item(70, 14).Value = 7000
$ws.Cells.Item(70, 15).Value = 7000
$ws.Cells.Item(70, 16).Value = 7000
$ws.Cells.Item(70, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(70, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(70, 19).Value = 700
$ws.Cells.Item(70, 20).Value = 10
